# Add team record (Wins/Losses/Ties) columns to the CLE_2023 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in row 1 (columns AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting (bold font + border) used by the other header cells
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Every player row (2-52) shares the same team record for the season
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 76   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 86   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
